$wb = $excel.ActiveWorkbook

$wsInvoice = $wb.Worksheets.Item("Historical Invoice Template")
$wsPO = $wb.Worksheets.Item("Historical PO Template")

# Historical Invoice Template (sheet1) - row 2
$wsInvoice.Range("B2").Value = "Invoice1735311"
$wsInvoice.Range("C2").Value = "Invoice1858875"
$wsInvoice.Range("E2").Value = "Aprilwi19"
$wsInvoice.Range("F2").Value = "Iris6cr3"

# Historical PO Template (sheet2) - row 2
$wsPO.Range("B2").Value = "Invoice1858875"
$wsPO.Range("C2").Value = "Invoice1735311"
$wsPO.Range("E2").Value = "Aprilwi19"
$wsPO.Range("F2").Value = "Iris6cr3"
